# Scheduled data refresh: update market-price / profit columns (H:N) on the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly
# pulled currentAveragePrice / LevePrice / LeveProfit figures.
#
# Columns: H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#          K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
#
# A few rows no longer have a computable NQ/HQ profit (the market data needed
# to compute them is unavailable), so their M/N cells are cleared entirely
# rather than zeroed; conversely a couple of rows gained a computable
# LeveProfitNQ value where previously there was none.

$wb = $excel.ActiveWorkbook

# Sheet index 1 (ALC)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(41, 8).Value = 1531.375
$ws.Cells.Item(41, 9).Value = 1846.4445
$ws.Cells.Item(41, 10).Value = 1342.3334
$ws.Cells.Item(41, 11).Value = 1846.4445
$ws.Cells.Item(41, 12).Value = 1342.3334
$ws.Cells.Item(41, 13).Value = -1406.4445
$ws.Cells.Item(41, 14).Value = -2222.3334
$ws.Cells.Item(53, 8).Value = 699.2857
$ws.Cells.Item(53, 9).Value = 320.66666
$ws.Cells.Item(53, 11).Value = 320.66666
$ws.Cells.Item(53, 13).Value = 316.33334
$ws.Cells.Item(62, 8).Value = 2699.8333
$ws.Cells.Item(62, 9).Value = 2699.8333
$ws.Cells.Item(62, 11).Value = 2699.8333
$ws.Cells.Item(62, 13).Value = -2075.8333
$ws.Cells.Item(65, 8).Value = 2699.8333
$ws.Cells.Item(65, 9).Value = 2699.8333
$ws.Cells.Item(65, 11).Value = 13499.1665
$ws.Cells.Item(65, 13).Value = -10379.1665
$ws.Cells.Item(70, 8).Value = 17540
$ws.Cells.Item(70, 10).Value = 39000
$ws.Cells.Item(70, 12).Value = 117000
$ws.Cells.Item(70, 14).Value = -117540
$ws.Cells.Item(73, 8).Value = 17540
$ws.Cells.Item(73, 10).Value = 39000
$ws.Cells.Item(73, 12).Value = 117000
$ws.Cells.Item(73, 14).Value = -118872
$ws.Cells.Item(138, 8).Value = 2577.1914
$ws.Cells.Item(138, 9).Value = 2407.4736
$ws.Cells.Item(138, 10).Value = 2692.3572
$ws.Cells.Item(138, 11).Value = 7222.4208
$ws.Cells.Item(138, 12).Value = 8077.071599999999
$ws.Cells.Item(138, 13).Value = -2082.4208
$ws.Cells.Item(138, 14).Value = -18357.0716
$ws.Cells.Item(141, 8).Value = 2058.5334
$ws.Cells.Item(141, 10).Value = 3000
$ws.Cells.Item(141, 12).Value = 9000
$ws.Cells.Item(141, 14).Value = -19360

# Sheet index 2 (ARM)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 360873.5
$ws.Cells.Item(2, 9).Value = 774384.0600000001
$ws.Cells.Item(2, 11).Value = 774384.0600000001
$ws.Cells.Item(2, 13).Value = -774271.0600000001
$ws.Cells.Item(61, 8).Value = 22731810
$ws.Cells.Item(61, 9).Value = 23813838
$ws.Cells.Item(61, 10).Value = 9223
$ws.Cells.Item(61, 11).Value = 23813838
$ws.Cells.Item(61, 12).Value = 9223
$ws.Cells.Item(61, 13).Value = -23813626
$ws.Cells.Item(61, 14).Value = -9647
$ws.Cells.Item(74, 8).Value = 38465950
$ws.Cells.Item(74, 9).Value = 41671200
$ws.Cells.Item(74, 11).Value = 41671200
$ws.Cells.Item(74, 13).Value = -41670326
$ws.Cells.Item(77, 8).Value = 38465950
$ws.Cells.Item(77, 9).Value = 41671200
$ws.Cells.Item(77, 11).Value = 208356000
$ws.Cells.Item(77, 13).Value = -208351632
$ws.Cells.Item(116, 8).Value = 360873.5
$ws.Cells.Item(116, 9).Value = 774384.0600000001
$ws.Cells.Item(116, 11).Value = 774384.0600000001
$ws.Cells.Item(116, 13).Value = -772090.0600000001
$ws.Cells.Item(136, 8).Value = 22731810
$ws.Cells.Item(136, 9).Value = 23813838
$ws.Cells.Item(136, 10).Value = 9223
$ws.Cells.Item(136, 11).Value = 71441514
$ws.Cells.Item(136, 12).Value = 27669
$ws.Cells.Item(136, 13).Value = -71438964
$ws.Cells.Item(136, 14).Value = -32769

# Sheet index 3 (BSM)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 360873.5
$ws.Cells.Item(3, 9).Value = 774384.0600000001
$ws.Cells.Item(3, 11).Value = 774384.0600000001
$ws.Cells.Item(3, 13).Value = -774270.0600000001
$ws.Cells.Item(86, 8).Value = 2255.2415
$ws.Cells.Item(86, 9).Value = 2375.55
$ws.Cells.Item(86, 10).Value = 1987.8889
$ws.Cells.Item(86, 11).Value = 2375.55
$ws.Cells.Item(86, 12).Value = 1987.8889
$ws.Cells.Item(86, 13).Value = -1252.55
$ws.Cells.Item(86, 14).Value = -4233.8889
$ws.Cells.Item(89, 8).Value = 2255.2415
$ws.Cells.Item(89, 9).Value = 2375.55
$ws.Cells.Item(89, 10).Value = 1987.8889
$ws.Cells.Item(89, 11).Value = 11877.75
$ws.Cells.Item(89, 12).Value = 9939.4445
$ws.Cells.Item(89, 13).Value = -6261.75
$ws.Cells.Item(89, 14).Value = -21171.4445
$ws.Cells.Item(105, 8).Value = 3347.84
$ws.Cells.Item(105, 9).Value = 3083.8333
$ws.Cells.Item(105, 10).Value = 4026.7144
$ws.Cells.Item(105, 11).Value = 3083.8333
$ws.Cells.Item(105, 12).Value = 4026.7144
$ws.Cells.Item(105, 13).Value = -1336.8333
$ws.Cells.Item(105, 14).Value = -7520.7144
$ws.Cells.Item(134, 8).Value = 15627568
$ws.Cells.Item(134, 9).Value = 16131038
$ws.Cells.Item(134, 11).Value = 48393114
$ws.Cells.Item(134, 13).Value = -48390579

# Sheet index 4 (CRP)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(8, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 8206.704
$ws.Cells.Item(31, 9).Value = 5571.8438
$ws.Cells.Item(31, 10).Value = 12039.228
$ws.Cells.Item(31, 11).Value = 5571.8438
$ws.Cells.Item(31, 12).Value = 12039.228
$ws.Cells.Item(31, 13).Value = -5276.8438
$ws.Cells.Item(31, 14).Value = -12629.228
$ws.Cells.Item(34, 8).Value = 8206.704
$ws.Cells.Item(34, 9).Value = 5571.8438
$ws.Cells.Item(34, 10).Value = 12039.228
$ws.Cells.Item(34, 11).Value = 5571.8438
$ws.Cells.Item(34, 12).Value = 12039.228
$ws.Cells.Item(34, 13).Value = -5369.8438
$ws.Cells.Item(34, 14).Value = -12443.228
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()
$ws.Cells.Item(58, 8).Value = 55570388
$ws.Cells.Item(58, 9).Value = 100021000
$ws.Cells.Item(58, 10).Value = 7124.75
$ws.Cells.Item(58, 11).Value = 100021000
$ws.Cells.Item(58, 12).Value = 7124.75
$ws.Cells.Item(58, 13).Value = -100020797
$ws.Cells.Item(58, 14).Value = -7530.75
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 6251142
$ws.Cells.Item(134, 9).Value = 7353942.5
$ws.Cells.Item(134, 11).Value = 22061827.5
$ws.Cells.Item(134, 13).Value = -22059292.5
$ws.Cells.Item(136, 8).Value = 55570388
$ws.Cells.Item(136, 9).Value = 100021000
$ws.Cells.Item(136, 10).Value = 7124.75
$ws.Cells.Item(136, 11).Value = 300063000
$ws.Cells.Item(136, 12).Value = 21374.25
$ws.Cells.Item(136, 13).Value = -300060450
$ws.Cells.Item(136, 14).Value = -26474.25

# Sheet index 5 (CUL)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(3, 8).Value = 199
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(14, 8).Value = 401.9
$ws.Cells.Item(14, 9).Value = 401.9
$ws.Cells.Item(14, 11).Value = 1205.7
$ws.Cells.Item(14, 13).Value = -1032.7
$ws.Cells.Item(18, 8).Value = 1681.3334
$ws.Cells.Item(18, 9).Value = 1022.25
$ws.Cells.Item(18, 11).Value = 3066.75
$ws.Cells.Item(18, 13).Value = -2897.75
$ws.Cells.Item(131, 8).Value = 1531.9166
$ws.Cells.Item(131, 9).Value = 1329.1875
$ws.Cells.Item(131, 10).Value = 1937.375
$ws.Cells.Item(131, 11).Value = 3987.5625
$ws.Cells.Item(131, 12).Value = 5812.125
$ws.Cells.Item(131, 13).Value = 1052.4375
$ws.Cells.Item(131, 14).Value = -15892.125

# Sheet index 6 (GSM)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 10106.941
$ws.Cells.Item(70, 9).Value = 9405.454
$ws.Cells.Item(70, 11).Value = 9405.454
$ws.Cells.Item(70, 13).Value = -9135.454
$ws.Cells.Item(73, 8).Value = 10106.941
$ws.Cells.Item(73, 9).Value = 9405.454
$ws.Cells.Item(73, 11).Value = 9405.454
$ws.Cells.Item(73, 13).Value = -8469.454
$ws.Cells.Item(80, 8).Value = 2445.6667
$ws.Cells.Item(80, 9).Value = 2335.0908
$ws.Cells.Item(80, 10).Value = 2749.75
$ws.Cells.Item(80, 11).Value = 2335.0908
$ws.Cells.Item(80, 12).Value = 2749.75
$ws.Cells.Item(80, 13).Value = -1337.0908
$ws.Cells.Item(80, 14).Value = -4745.75
$ws.Cells.Item(83, 8).Value = 2445.6667
$ws.Cells.Item(83, 9).Value = 2335.0908
$ws.Cells.Item(83, 10).Value = 2749.75
$ws.Cells.Item(83, 11).Value = 11675.454
$ws.Cells.Item(83, 12).Value = 13748.75
$ws.Cells.Item(83, 13).Value = -6683.454
$ws.Cells.Item(83, 14).Value = -23732.75
$ws.Cells.Item(122, 8).Value = 3858
$ws.Cells.Item(122, 9).Value = 560.6667
$ws.Cells.Item(122, 11).Value = 1682.0001
$ws.Cells.Item(122, 13).Value = 767.9999

# Sheet index 7 (LTW)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 1751.8529
$ws.Cells.Item(16, 9).Value = 484.14285
$ws.Cells.Item(16, 11).Value = 484.14285
$ws.Cells.Item(16, 13).Value = -314.14285
$ws.Cells.Item(22, 8).Value = 3642.7144
$ws.Cells.Item(22, 9).Value = 2899.8
$ws.Cells.Item(22, 11).Value = 2899.8
$ws.Cells.Item(22, 13).Value = -2604.8
$ws.Cells.Item(27, 8).Value = 3642.7144
$ws.Cells.Item(27, 9).Value = 2899.8
$ws.Cells.Item(27, 11).Value = 2899.8
$ws.Cells.Item(27, 13).Value = -2792.8
$ws.Cells.Item(40, 8).Value = 6331
$ws.Cells.Item(40, 9).Value = 3500
$ws.Cells.Item(40, 11).Value = 3500
$ws.Cells.Item(40, 13).Value = -3364
$ws.Cells.Item(55, 8).Value = 189.5625
$ws.Cells.Item(55, 9).Value = 132.61111
$ws.Cells.Item(55, 10).Value = 262.7857
$ws.Cells.Item(55, 11).Value = 132.61111
$ws.Cells.Item(55, 12).Value = 262.7857
$ws.Cells.Item(55, 13).Value = 40.38889
$ws.Cells.Item(55, 14).Value = -608.7857
$ws.Cells.Item(64, 8).Value = 20331.334
$ws.Cells.Item(64, 9).Value = 17999.5
$ws.Cells.Item(64, 10).Value = 24995
$ws.Cells.Item(64, 11).Value = 17999.5
$ws.Cells.Item(64, 12).Value = 24995
$ws.Cells.Item(64, 13).Value = -17774.5
$ws.Cells.Item(64, 14).Value = -25445
$ws.Cells.Item(67, 8).Value = 20331.334
$ws.Cells.Item(67, 9).Value = 17999.5
$ws.Cells.Item(67, 10).Value = 24995
$ws.Cells.Item(67, 11).Value = 17999.5
$ws.Cells.Item(67, 12).Value = 24995
$ws.Cells.Item(67, 13).Value = -17219.5
$ws.Cells.Item(67, 14).Value = -26555
$ws.Cells.Item(82, 8).Value = 1155
$ws.Cells.Item(82, 9).Value = 1147
$ws.Cells.Item(82, 11).Value = 1147
$ws.Cells.Item(82, 13).Value = -786
$ws.Cells.Item(85, 8).Value = 1155
$ws.Cells.Item(85, 9).Value = 1147
$ws.Cells.Item(85, 11).Value = 1147
$ws.Cells.Item(85, 13).Value = 101
$ws.Cells.Item(93, 8).Value = 1283.3334
$ws.Cells.Item(93, 9).Value = 425
$ws.Cells.Item(93, 11).Value = 425
$ws.Cells.Item(93, 13).Value = 823
$ws.Cells.Item(132, 8).Value = 25275726
$ws.Cells.Item(132, 9).Value = 30014176
$ws.Cells.Item(132, 10).Value = 3996.6667
$ws.Cells.Item(132, 11).Value = 90042528
$ws.Cells.Item(132, 12).Value = 11990.0001
$ws.Cells.Item(132, 13).Value = -90039998
$ws.Cells.Item(132, 14).Value = -17050.0001
$ws.Cells.Item(136, 8).Value = 2920.7222
$ws.Cells.Item(136, 9).Value = 1090
$ws.Cells.Item(136, 11).Value = 3270
$ws.Cells.Item(136, 13).Value = -720

# Sheet index 8 (WVR)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 2726.3076
$ws.Cells.Item(81, 9).Value = 2726.3076
$ws.Cells.Item(81, 11).Value = 5452.6152
$ws.Cells.Item(81, 13).Value = -4391.6152
$ws.Cells.Item(84, 8).Value = 2726.3076
$ws.Cells.Item(84, 9).Value = 2726.3076
$ws.Cells.Item(84, 11).Value = 27263.076
$ws.Cells.Item(84, 13).Value = -21959.076
